$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a "daily auto push" time-series log. Two more rows for
# 2026/02/25 (times 19:00 and 22:00) need to be inserted right after the
# existing 2026/02/25 16:00 entry (old row 882) and before the 2026/12/29
# block (old row 883). Every row from the old row 883 onward shifts down
# by 2 (old 924 -> new 926).

$ws.Range("A883:A884").EntireRow.Insert(-4121)

# New row 883: 2026/02/25, 水, 19:00, ranking 201
$ws.Cells.Item(883, 1).NumberFormat = "@"
$ws.Cells.Item(883, 1).Value = "2026/02/25"
$ws.Cells.Item(883, 1).Style = "Normal"
$ws.Cells.Item(883, 2).Value = "水"
$ws.Cells.Item(883, 3).Value = 19
$ws.Cells.Item(883, 4).Value = 201

# New row 884: 2026/02/25, 水, 22:00, ranking 201
$ws.Cells.Item(884, 1).NumberFormat = "@"
$ws.Cells.Item(884, 1).Value = "2026/02/25"
$ws.Cells.Item(884, 1).Style = "Normal"
$ws.Cells.Item(884, 2).Value = "水"
$ws.Cells.Item(884, 3).Value = 22
$ws.Cells.Item(884, 4).Value = 201

Write-Output "done"
